$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.927.17'
$ws.Range('E2').Value = '  -5.46%  '
$ws.Range('D3').Value = '3.165.86'
$ws.Range('E3').Value = '  -7.03%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').Value = '513.84'
$ws.Range('E5').Value = '  -3.47%  '
$ws.Range('D6').Value = '168.35'
$ws.Range('E6').Value = '  -9.71%  '
$ws.Range('D7').Value = '0.583'
$ws.Range('E7').Value = '  -4.65%  '
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').Value = '3.169.64'
$ws.Range('E9').Value = '  -6.91%  '
$ws.Range('D10').Value = '0.588'
$ws.Range('E10').Value = '  -7.06%  '
$ws.Range('D11').Value = '51.53'
$ws.Range('E11').Value = '  -12.74%  '
$ws.Range('D12').Value = '0.127'
$ws.Range('E12').Value = '  -5.80%  '
$ws.Range('D13').Value = '0.0000247'
$ws.Range('E13').Value = '  -4.35%  '
$ws.Range('D14').Value = '8.71'
$ws.Range('E14').Value = '  -6.83%  '
$ws.Range('D15').Value = '3.698.28'
$ws.Range('E15').Value = '  -6.14%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.183.11'
$ws.Range('E16').Value = '  -6.47%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '0.113'
$ws.Range('E17').Value = '  -8.30%  '
$ws.Range('D18').Value = '62.026.26'
$ws.Range('E18').Value = '  -4.80%  '
$ws.Range('D19').Value = '16.83'
$ws.Range('E19').Value = '  -4.70%  '
$ws.Range('D20').Value = '10.71'
$ws.Range('E20').Value = '  -5.37%  '
$ws.Range('D21').Value = '0.940'
$ws.Range('E21').Value = '  -4.32%  '
$ws.Range('D22').Value = '358.06'
$ws.Range('E22').Value = '  -4.79%  '
$ws.Range('D23').Value = '3.66'
$ws.Range('E23').Value = '  -3.25%  '
$ws.Range('D24').Value = '79.29'
$ws.Range('E24').Value = '  -3.90%  '
$ws.Range('D25').Value = '10.77'
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('D26').Value = '6.11'
$ws.Range('E26').Value = '  +4.35%  '
$ws.Range('D27').Value = '3.78'
$ws.Range('E27').Value = '  +1.67%  '
$ws.Range('D28').Value = '2.55'
$ws.Range('E28').Value = '  -5.18%  '
$ws.Range('D29').Value = '10.91'
$ws.Range('E29').Value = '  -7.59%  '
$ws.Range('D30').Value = '7.99'
$ws.Range('E30').Value = '  -7.42%  '
$ws.Range('D31').Value = '637.50'
$ws.Range('E31').Value = '  -7.13%  '
$ws.Range('D32').Value = '27.72'
$ws.Range('E32').Value = '  -7.40%  '
$ws.Range('D33').Value = '6.26'
$ws.Range('E33').Value = '  -8.47%  '
$ws.Range('D34').Value = '10.99'
$ws.Range('E34').Value = '  -2.96%  '
$ws.Range('D35').Value = '0.102'
$ws.Range('E35').Value = '  -4.88%  '
$ws.Range('D36').Value = '56.55'
$ws.Range('E36').Value = '  -8.44%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').Value = '35.58'
$ws.Range('E38').Value = '  -3.41%  '
$ws.Range('D39').Value = '0.364'
$ws.Range('E39').Value = '  -6.20%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').Value = '0.0″0679'
$ws.Range('E41').Value = '  +7.41%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.120'
$ws.Range('E42').Value = '  -6.78%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.830.18'
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('D44').Value = '2.43'
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('E45').Value = '  -3.24%  '
$ws.Range('D46').Value = '0.0378'
$ws.Range('E46').Value = '  -5.86%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '2.72'
$ws.Range('E47').Value = '  +2.18%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = '2.49'
$ws.Range('E48').Value = '  -11.36%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '2.93'
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('D50').Value = '133.31'
$ws.Range('E50').Value = '  -3.15%  '
$ws.Range('D51').Value = '0.120'
$ws.Range('E51').Value = '  -5.18%  '
